$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Renumber the scenario labels in column A (rows 2-7) down by one:
# b1 -> b0, b2 -> b1, b3 -> b2, b4 -> b3, b5 -> b4, b6 -> b5
$ws.Range("A2").Value = "b0"
$ws.Range("A3").Value = "b1"
$ws.Range("A4").Value = "b2"
$ws.Range("A5").Value = "b3"
$ws.Range("A6").Value = "b4"
$ws.Range("A7").Value = "b5"
